$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) data, and the
# swapped MXToken/RenderToken rows (38-39), per the latest scrape.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.910.50'
$ws.Range("E2").Value = '  +1.70%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.882.88'
$ws.Range("E3").Value = '  +2.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.73%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.79'
$ws.Range("E5").Value = '  -1.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4933'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.30'
$ws.Range("E8").Value = '  -0.61%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2902'
$ws.Range("E9").Value = '  +3.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06587'
$ws.Range("E10").Value = '  +2.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.891.43'
$ws.Range("E11").Value = '  +2.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.83'
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07187'
$ws.Range("E13").Value = '  +1.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6673'
$ws.Range("E14").Value = '  +2.31%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '85.24'
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.797'
$ws.Range("E16").Value = '  +1.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.975.82'
$ws.Range("E17").Value = '  +1.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007833'
$ws.Range("E18").Value = '  +6.62%  '

$ws.Range("E19").Value = '  +0.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.73'
$ws.Range("E20").Value = '  +2.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.133.89'
$ws.Range("E21").Value = '  +3.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.62%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.733'
$ws.Range("E23").Value = '  +3.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.561'
$ws.Range("E24").Value = '  +2.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.076'
$ws.Range("E25").Value = '  +2.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '148.41'
$ws.Range("E26").Value = '  +3.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '133.42'
$ws.Range("E27").Value = '  +1.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.65'
$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.922'
$ws.Range("E29").Value = '  +0.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.378'
$ws.Range("E30").Value = '  -1.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.159'
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08663'
$ws.Range("E32").Value = '  +3.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.915'
$ws.Range("E33").Value = '  +2.94%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05089'
$ws.Range("E34").Value = '  +2.93%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.107'
$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7017'
$ws.Range("E36").Value = '  +3.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.683'
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.691'
$ws.Range("E38").Value = '  -0.90%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.198'
$ws.Range("E39").Value = '  -3.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9334'
$ws.Range("E40").Value = '  -2.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01640'
$ws.Range("E41").Value = '  +2.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.045'
$ws.Range("E42").Value = '  -2.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9990'
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.50'
$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4161'
$ws.Range("E45").Value = '  +1.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.430'
$ws.Range("E46").Value = '  +2.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1256'
$ws.Range("E47").Value = '  +2.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05719'
$ws.Range("E48").Value = '  +2.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.48'
$ws.Range("E49").Value = '  +1.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.218'
$ws.Range("E50").Value = '  +1.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3693'
$ws.Range("E51").Value = '  +1.81%  '
